$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 48: date label (column A) + gold-price summary (column B), appended
# after the existing last row (47), matching the sheet's existing text-only
# convention for the date column.
#
# A plain Value/Value2 assignment of "01-11-2025" gets auto-detected by the
# COM layer as a date (mm-dd-yyyy) and rewritten to a date serial number with
# a brand-new number-format style, which would diverge from the target sheet
# (which keeps date labels as literal shared strings with the existing
# inherited column style). Writing it as a literal-text formula first keeps
# it textual, then collapsing the formula to its cached value via a
# values-only paste removes the formula while preserving the original style.
$cellA = $ws.Cells.Item(48, 1)
$cellA.Formula = "=""01-11-2025"""
$cellA.Copy()
$cellA.PasteSpecial(-4163)  # xlPasteValues

$ws.Cells.Item(48, 2).Value = "The price of gold in India today is ₹12,300 per gram for 24 karat gold, ₹11,275 per gram for 22 karat gold and ₹9,225 per gram for 18 karat gold (also called 999 gold)."
